# conditionRuntimeException-template.docx used to encode its `m:if` / `m:endif`
# markers as real Word fields (w:fldChar begin/end wrapping w:instrText runs).
# The parser was switched to TokenIteratorFieldRewriterSplit, which instead
# expects the markers as plain literal text of the form "{...}" split across
# runs that mirror the field code's original run layout. This script rewrites
# both fields in place, run-by-run, while leaving every other paragraph
# (including paragraph properties / bookmarks) untouched.

$d = $word.ActiveDocument

function ConvertTo-XmlText($s) {
    # Minimal XML-escape for text going back into a w:t element.
    return $s.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
}

# Builds the "<w:r><w:t>...</w:t></w:r>" sequence for one field: $tokens holds
# the literal text that used to live in each original w:instrText run (in
# document order) and $preserveFlags says whether that original run carried
# xml:space="preserve". The first token gets "{" prepended and the last gets
# "}" appended -- exactly what the new TokenIteratorFieldRewriterSplit output
# looks like -- while the preserve flags stay tied to the *original* run they
# came from (not recomputed from the brace-augmented text).
function Build-FieldRunsXml($tokens, $preserveFlags) {
    $n = $tokens.Count
    $xml = ""
    for ($i = 0; $i -lt $n; $i++) {
        $tok = $tokens[$i]
        if ($i -eq 0) { $tok = "{" + $tok }
        if ($i -eq ($n - 1)) { $tok = $tok + "}" }
        $escaped = ConvertTo-XmlText $tok
        if ($preserveFlags[$i]) {
            $xml += "<w:r><w:t xml:space=`"preserve`">$escaped</w:t></w:r>"
        } else {
            $xml += "<w:r><w:t>$escaped</w:t></w:r>"
        }
    }
    return $xml
}

# Replaces the field $f with literal runs built from $tokens/$preserveFlags:
# the new runs are inserted as a zero-length InsertXML right at the field's
# begin mark (this leaves the owning paragraph's w:pPr/bookmarks untouched),
# then the original field (begin fldChar + instrText run(s) + end fldChar) is
# removed with Field.Delete().
function Convert-FieldToLiteralText($f, $tokens, $preserveFlags) {
    $runsXml = Build-FieldRunsXml $tokens $preserveFlags
    $insPos = $f.Code.Start - 1
    $insRange = $d.Range($insPos, $insPos)
    $packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' + $runsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $insRange.InsertXML($packageXml)
    $f.Delete()
}

# Walk the document's fields and rewrite the two recognised M2Doc markers.
# Fields.Item(1) is re-fetched every time because Field.Delete() shifts the
# collection down by one, so the next unseen field always resettles at index 1.
$guard = 0
while ($d.Fields.Count -gt 0 -and $guard -lt 50) {
    $guard = $guard + 1
    $f = $d.Fields.Item(1)
    $code = $f.Code.Text.Trim()

    if ($code -eq "m:if 1/0 = 42") {
        $tokens = @("m:if ", "1/0", " ", "=", " ", "42")
        $preserve = @($true, $false, $true, $false, $true, $false)
        Convert-FieldToLiteralText $f $tokens $preserve
    } elseif ($code -eq "m:endif") {
        $tokens = @("m:endif")
        $preserve = @($true)
        Convert-FieldToLiteralText $f $tokens $preserve
    } else {
        # Unrecognised field: leave it alone and stop touching the Fields
        # collection so we don't loop forever on something we don't understand.
        break
    }
}
